$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.165.90'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +6.96%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.619.15'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.73%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '416.19'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.94'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.661'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +4.81%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.604.47'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.49%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.760'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.68%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +24.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000425'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +89.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '42.08'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.83'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.194.65'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.99%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.606.54'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.98'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.74%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.077.95'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +7.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.32'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '457.83'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '89.16'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.29'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -6.53%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.28'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '35.37'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +5.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.96'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.81%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.19'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.27%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.06%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '40.49'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.06%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -7.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.997'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '55.79'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0786'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +17.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0484'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.84%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +8.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '148.59'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.95'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.23'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.82%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.58'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +10.23%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.24'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -7.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.169'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +20.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.300'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.35%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.94'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.39%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +13.09%  '
